$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Rows.Item(1).Cells.Item(1).Range.Text  = "97÷3=32, 1"
$t.Rows.Item(1).Cells.Item(2).Range.Text  = "39÷3=13, 0"
$t.Rows.Item(1).Cells.Item(3).Range.Text  = "10÷2=5, 0"
$t.Rows.Item(1).Cells.Item(4).Range.Text  = "73÷4=18, 1"
$t.Rows.Item(1).Cells.Item(5).Range.Text  = "68÷3=22, 2"

$t.Rows.Item(5).Cells.Item(1).Range.Text  = "38÷6=6, 2"
$t.Rows.Item(5).Cells.Item(2).Range.Text  = "24÷3=8, 0"
$t.Rows.Item(5).Cells.Item(3).Range.Text  = "57÷4=14, 1"
$t.Rows.Item(5).Cells.Item(4).Range.Text  = "51÷7=7, 2"
$t.Rows.Item(5).Cells.Item(5).Range.Text  = "76÷3=25, 1"

$t.Rows.Item(9).Cells.Item(1).Range.Text  = "38÷7=5, 3"
$t.Rows.Item(9).Cells.Item(2).Range.Text  = "89÷2=44, 1"
$t.Rows.Item(9).Cells.Item(3).Range.Text  = "47÷3=15, 2"
$t.Rows.Item(9).Cells.Item(4).Range.Text  = "54÷9=6, 0"
$t.Rows.Item(9).Cells.Item(5).Range.Text  = "93÷4=23, 1"

$t.Rows.Item(13).Cells.Item(1).Range.Text = "22÷6=3, 4"
$t.Rows.Item(13).Cells.Item(2).Range.Text = "85÷8=10, 5"
$t.Rows.Item(13).Cells.Item(3).Range.Text = "28÷4=7, 0"
$t.Rows.Item(13).Cells.Item(4).Range.Text = "73÷2=36, 1"
$t.Rows.Item(13).Cells.Item(5).Range.Text = "97÷5=19, 2"

$t.Rows.Item(17).Cells.Item(1).Range.Text = "84÷9=9, 3"
$t.Rows.Item(17).Cells.Item(2).Range.Text = "88÷3=29, 1"
$t.Rows.Item(17).Cells.Item(3).Range.Text = "86÷5=17, 1"
$t.Rows.Item(17).Cells.Item(4).Range.Text = "36÷5=7, 1"
$t.Rows.Item(17).Cells.Item(5).Range.Text = "29÷5=5, 4"
